$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.012617333333333
$ws.Range("H2").Value = 3.037852
$ws.Range("I2").Value = 0.0186050446061446
$ws.Range("J2").Value = 0.0186050446061446
$ws.Range("M2").Value = 2.134850333333333
$ws.Range("N2").Value = 6.404551
$ws.Range("O2").Value = 0.03148693319884856
$ws.Range("P2").Value = 0.03148693319884856
$ws.Range("Q2").Value = 2.161786451605777
$ws.Range("R2").Value = 19.456078064452
$ws.Range("S2").Value = 0.0005858157966752727
$ws.Range("T2").Value = 0.0005858157966752729
$ws.Range("G3").Value = 1.012617333333333
$ws.Range("H3").Value = 3.037852
$ws.Range("I3").Value = 0.0186050446061446
$ws.Range("J3").Value = 0.0186050446061446
$ws.Range("O3").Value = 0.004432963048734841
$ws.Range("P3").Value = 0.004432963048734842
$ws.Range("Q3").Value = 0.3043522657066666
$ws.Range("R3").Value = 2.73917039136
$ws.Range("S3").Value = 0.00008247547525910247
$ws.Range("T3").Value = 0.00008247547525910248
$ws.Range("G4").Value = 1.012617333333333
$ws.Range("H4").Value = 3.037852
$ws.Range("I4").Value = 0.0186050446061446
$ws.Range("J4").Value = 0.0186050446061446
$ws.Range("M4").Value = 1.073938666666667
$ws.Range("N4").Value = 3.221816
$ws.Range("O4").Value = 0.01583953428912994
$ws.Range("P4").Value = 0.01583953428912995
$ws.Range("Q4").Value = 1.087488908803555
$ws.Range("R4").Value = 9.787400179232
$ws.Range("S4").Value = 0.0002946952419898196
$ws.Range("T4").Value = 0.0002946952419898196
$ws.Range("G5").Value = 1.012617333333333
$ws.Range("H5").Value = 3.037852
$ws.Range("I5").Value = 0.0186050446061446
$ws.Range("J5").Value = 0.0186050446061446
$ws.Range("M5").Value = 64.29180266666667
$ws.Range("N5").Value = 192.875408
$ws.Range("O5").Value = 0.9482405694632866
$ws.Range("P5").Value = 0.9482405694632867
$ws.Range("Q5").Value = 65.10299377151289
$ws.Range("R5").Value = 585.9269439436159
$ws.Range("S5").Value = 0.0176420580922204
$ws.Range("T5").Value = 0.01764205809222041
$ws.Range("I6").Value = 0.7824865355506074
$ws.Range("J6").Value = 0.7824865355506075
$ws.Range("M6").Value = 2.134850333333333
$ws.Range("N6").Value = 6.404551
$ws.Range("O6").Value = 0.03148693319884856
$ws.Range("P6").Value = 0.03148693319884856
$ws.Range("Q6").Value = 90.91989978667287
$ws.Range("R6").Value = 818.2790980800559
$ws.Range("S6").Value = 0.02463810127388041
$ws.Range("T6").Value = 0.02463810127388042
$ws.Range("I7").Value = 0.7824865355506074
$ws.Range("J7").Value = 0.7824865355506075
$ws.Range("O7").Value = 0.004432963048734841
$ws.Range("P7").Value = 0.004432963048734842
$ws.Range("S7").Value = 0.003468733898228384
$ws.Range("T7").Value = 0.003468733898228385
$ws.Range("I8").Value = 0.7824865355506074
$ws.Range("J8").Value = 0.7824865355506075
$ws.Range("M8").Value = 1.073938666666667
$ws.Range("N8").Value = 3.221816
$ws.Range("O8").Value = 0.01583953428912994
$ws.Range("P8").Value = 0.01583953428912995
$ws.Range("Q8").Value = 45.73734955832177
$ws.Range("R8").Value = 411.636146024896
$ws.Range("S8").Value = 0.01239422231063634
$ws.Range("T8").Value = 0.01239422231063635
$ws.Range("I9").Value = 0.7824865355506074
$ws.Range("J9").Value = 0.7824865355506075
$ws.Range("M9").Value = 64.29180266666667
$ws.Range("N9").Value = 192.875408
$ws.Range("O9").Value = 0.9482405694632866
$ws.Range("P9").Value = 0.9482405694632867
$ws.Range("Q9").Value = 2738.086208802716
$ws.Range("R9").Value = 24642.77587922445
$ws.Range("S9").Value = 0.7419854780678623
$ws.Range("T9").Value = 0.7419854780678624
$ws.Range("G10").Value = 10.82599466666667
$ws.Range("H10").Value = 32.477984
$ws.Range("I10").Value = 0.198908419843248
$ws.Range("J10").Value = 0.198908419843248
$ws.Range("M10").Value = 2.134850333333333
$ws.Range("N10").Value = 6.404551
$ws.Range("O10").Value = 0.03148693319884856
$ws.Range("P10").Value = 0.03148693319884856
$ws.Range("Q10").Value = 23.11187832279822
$ws.Range("R10").Value = 208.006904905184
$ws.Range("S10").Value = 0.006263016128292873
$ws.Range("T10").Value = 0.006263016128292874
$ws.Range("G11").Value = 10.82599466666667
$ws.Range("H11").Value = 32.477984
$ws.Range("I11").Value = 0.198908419843248
$ws.Range("J11").Value = 0.198908419843248
$ws.Range("O11").Value = 0.004432963048734841
$ws.Range("P11").Value = 0.004432963048734842
$ws.Range("Q11").Value = 3.253860957013333
$ws.Range("R11").Value = 29.28474861312
$ws.Range("S11").Value = 0.0008817536752473543
$ws.Range("T11").Value = 0.0008817536752473545
$ws.Range("G12").Value = 10.82599466666667
$ws.Range("H12").Value = 32.477984
$ws.Range("I12").Value = 0.198908419843248
$ws.Range("J12").Value = 0.198908419843248
$ws.Range("M12").Value = 1.073938666666667
$ws.Range("N12").Value = 3.221816
$ws.Range("O12").Value = 0.01583953428912994
$ws.Range("P12").Value = 0.01583953428912995
$ws.Range("Q12").Value = 11.62645427766044
$ws.Range("R12").Value = 104.638088498944
$ws.Range("S12").Value = 0.003150616736503782
$ws.Range("T12").Value = 0.003150616736503783
$ws.Range("G13").Value = 10.82599466666667
$ws.Range("H13").Value = 32.477984
$ws.Range("I13").Value = 0.198908419843248
$ws.Range("J13").Value = 0.198908419843248
$ws.Range("M13").Value = 64.29180266666667
$ws.Range("N13").Value = 192.875408
$ws.Range("O13").Value = 0.9482405694632866
$ws.Range("P13").Value = 0.9482405694632867
$ws.Range("Q13").Value = 696.0227127797191
$ws.Range("R13").Value = 6264.204415017472
$ws.Range("S13").Value = 0.188613033303204
$ws.Range("T13").Value = 0.188613033303204
